$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh scraped price / 1h-volume figures (rows 2-47)
$ws.Range("D2").Value = "29.133.80"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "1.838.23"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "239.96"
$ws.Range("E5").Value = "  -2.31%  "

$ws.Range("D6").Value = "0.6809"
$ws.Range("E6").Value = "  -2.85%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.2987"
$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").Value = "'0.07440"
$ws.Range("E9").Value = "  -3.98%  "

$ws.Range("D10").Value = "23.17"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").Value = "0.07645"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").Value = "1.834.59"
$ws.Range("E12").Value = "  -1.67%  "

$ws.Range("D13").Value = "5.025"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D14").Value = "0.6801"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "86.99"
$ws.Range("E15").Value = "  -5.94%  "

$ws.Range("D16").Value = "6.146"
$ws.Range("E16").Value = "  -6.92%  "

$ws.Range("D17").Value = "29.145.49"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").Value = "2.088.34"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "229.77"
$ws.Range("E20").Value = "  -5.26%  "

$ws.Range("D21").Value = "12.49"
$ws.Range("E21").Value = "  -2.30%  "

$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "7.342"
$ws.Range("E23").Value = "  -3.87%  "

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "161.57"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").Value = "0.1428"
$ws.Range("E26").Value = "  -5.67%  "

$ws.Range("D27").Value = "8.704"
$ws.Range("E27").Value = "  -2.54%  "

$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  -2.52%  "

$ws.Range("D30").Value = "4.247"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("D31").Value = "4.143"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "0.05339"
$ws.Range("E33").Value = "  +4.51%  "

$ws.Range("D34").Value = "0.7523"
$ws.Range("E34").Value = "  -4.31%  "

$ws.Range("D35").Value = "1.843"
$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("D36").Value = "'1.130"
$ws.Range("E36").Value = "  -2.44%  "

$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "1.312.81"
$ws.Range("E38").Value = "  -1.53%  "

$ws.Range("D39").Value = "'0.01820"
$ws.Range("E39").Value = "  -3.25%  "

$ws.Range("D40").Value = "2.719"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").Value = "0.9366"
$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("D42").Value = "6.072"
$ws.Range("E42").Value = "  +1.21%  "

$ws.Range("D43").Value = "105.32"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "0.08223"
$ws.Range("E45").Value = "  +30.18%  "

$ws.Range("D46").Value = "1.981.69"
$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("D47").Value = "0.5177"
$ws.Range("E47").Value = "  -0.85%  "

# Coin-ranking moved: Aave/RenderToken swapped, Cronos dropped for
# BabyDogeCoin, EnergySwap shifted down a row (rows 48-51)
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.774"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "64.16"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000121"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.436"
$ws.Range("E51").Value = "  -3.69%  "

